$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2024-02-20 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-21 Wednesday", 2) | Out-Null

$tbl = $d.Tables.Item(1)

# Row 1
$tbl.Cell(1, 1).Range.Text = "93÷6=15, 3"
$tbl.Cell(1, 2).Range.Text = "45÷4=11, 1"
$tbl.Cell(1, 3).Range.Text = "11÷2=5, 1"
$tbl.Cell(1, 4).Range.Text = "45÷8=5, 5"
$tbl.Cell(1, 5).Range.Text = "36÷8=4, 4"

# Row 5
$tbl.Cell(5, 1).Range.Text = "96÷7=13, 5"
$tbl.Cell(5, 2).Range.Text = "43÷6=7, 1"
$tbl.Cell(5, 3).Range.Text = "83÷9=9, 2"
$tbl.Cell(5, 4).Range.Text = "93÷2=46, 1"
$tbl.Cell(5, 5).Range.Text = "26÷4=6, 2"

# Row 9
$tbl.Cell(9, 1).Range.Text = "63÷4=15, 3"
$tbl.Cell(9, 2).Range.Text = "36÷2=18, 0"
$tbl.Cell(9, 3).Range.Text = "57÷6=9, 3"
$tbl.Cell(9, 4).Range.Text = "50÷7=7, 1"
$tbl.Cell(9, 5).Range.Text = "70÷8=8, 6"

# Row 13
$tbl.Cell(13, 1).Range.Text = "26÷7=3, 5"
$tbl.Cell(13, 2).Range.Text = "92÷5=18, 2"
$tbl.Cell(13, 3).Range.Text = "15÷3=5, 0"
$tbl.Cell(13, 4).Range.Text = "39÷4=9, 3"
$tbl.Cell(13, 5).Range.Text = "62÷9=6, 8"

# Row 17
$tbl.Cell(17, 1).Range.Text = "46÷8=5, 6"
$tbl.Cell(17, 2).Range.Text = "91÷5=18, 1"
$tbl.Cell(17, 3).Range.Text = "33÷4=8, 1"
$tbl.Cell(17, 4).Range.Text = "11÷4=2, 3"
$tbl.Cell(17, 5).Range.Text = "27÷2=13, 1"
